$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.039.63"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.660.43"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.08%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "514.31"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.30"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("E8").Value = "  +1.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.658.08"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.32"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.106"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.337"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("E13").Value = "  -1.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.122.82"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.000.81"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.09"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.77%  "
$ws.Range("E17").Value = "  +0.95%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.651.17"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.29%  "
$ws.Range("E19").Value = "  -0.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "344.08"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.30%  "
$ws.Range("E21").Value = "  +0.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.10"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.93%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.05"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("E25").Value = "  +1.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.740.45"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  -0.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.160"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0808"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.12"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.98%  "
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.43"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +9.12%  "
$ws.Range("E33").Value = "  +0.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.90"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "149.97"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.01"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +12.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.04"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +4.05%  "
$ws.Range("E38").Value = "  +3.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.856"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.51"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.69"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.51%  "
$ws.Range("E42").Value = "  +0.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "281.64"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.56%  "
$ws.Range("E44").Value = "  -1.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.998"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0985"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.41%  "
$ws.Range("E47").Value = "  +2.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0536"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.59%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.71"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.42%  "
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.28"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.03%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0230"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.10%  "
